# Actualizacion 2 de abril de 2024 - Lap HP
# Se actualiza el repositorio con el material del curso.
#
# Corrige la columna "Version C" de la hoja de respuestas (Tabla 1).
# La tabla tiene 21 filas: fila 1 es el encabezado (Pregunta, Version A,
# Version B, Version C); la fila N+1 contiene la respuesta de la
# pregunta N. La columna 4 es "Version C".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$answers = @{
    2  = "D"
    3  = "D"
    4  = "D"
    5  = "D"
    6  = "B"
    7  = "A"
    8  = "B"
    9  = "A"
    10 = "C"
    12 = "C"
    13 = "A"
    15 = "A"
    16 = "C"
    17 = "D"
    18 = "B"
    19 = "C"
    20 = "B"
}

foreach ($q in $answers.Keys) {
    $row = $q + 1
    $cell = $t.Cell($row, 4)
    $cell.Range.Text = $answers[$q]
}
